$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Prescaler value (B6): 4999 -> 9
$ws.Range("B6").Value = 9

# Update Time base Required value (B14): 1 -> 0.00001 (1E-5)
$ws.Range("B14").Value = 0.00001

# Move the active selection to B12
$ws.Activate()
$ws.Range("B12").Select()
